# Workbook / sheet handles
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rfid_item")

# Shift the "item" column values down one position for the existing rfid rows
# (endpoint accepts item param / position adjustment)
$ws.Range("B5").Value  = "PXmYk7IzzsrHFMq5j70o"
$ws.Range("B6").Value  = "RMWLUuACH72OuqSPYQDk"
$ws.Range("B7").Value  = "VfgrHcX6LvHuAvkJtdgU"
$ws.Range("B8").Value  = "YvxptylcQC7o6s7fK7H9"
$ws.Range("B9").Value  = "oZGiQLJMymfo2Mc4KJYm"

# New mapping row for an additional rfid/item pair
$ws.Range("A10").Value = "EEEEEEEE"
$ws.Range("B10").Value = "rxRod7cigQjBK9dDmlHv"

# Grow the table (Table2) so it includes the newly added row
$lo = $ws.ListObjects.Item("Table2")
$lo.Resize($ws.Range("A1:B10"))

# Update the current selection on the sheet
$ws.Range("B14").Select()
